$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.91
$ws.Range("AF3").Value = 15

# Row 4
$ws.Range("J4").Value = 2.38
$ws.Range("K4").Value = 2.2
$ws.Range("L4").Value = 5
$ws.Range("Y4").Value = 1.4
$ws.Range("Z4").Value = 2.75
$ws.Range("AA4").Value = 2.1
$ws.Range("AB4").Value = 1.67
$ws.Range("AD4").Value = 7
$ws.Range("AI4").Value = 8
$ws.Range("AP4").Value = 17
$ws.Range("AS4").Value = 51

# Row 5
$ws.Range("I5").Value = 1.82
$ws.Range("J5").Value = 4.5
$ws.Range("T5").Value = 1.9
$ws.Range("W5").Value = 2.62
$ws.Range("X5").Value = 1.37
$ws.Range("AA5").Value = 1.6
$ws.Range("AB5").Value = 2.05
$ws.Range("AJ5").Value = 6.5
$ws.Range("AN5").Value = 7.7
$ws.Range("AO5").Value = 9.25
$ws.Range("AQ5").Value = 16
$ws.Range("AS5").Value = 22
